$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.88321066666667
$ws.Range("H2").Value = 107.649632
$ws.Range("I2").Value = 0.08317795499144418
$ws.Range("J2").Value = 0.08448843719082051
$ws.Range("M2").Value = 22.323506
$ws.Range("N2").Value = 66.970518
$ws.Range("O2").Value = 0.2680681429232296
$ws.Range("P2").Value = 0.2757701347495435
$ws.Range("Q2").Value = 801.0390686165973
$ws.Range("R2").Value = 7209.351617549377
$ws.Range("S2").Value = 0.02229735992670841
$ws.Range("T2").Value = 0.02329938770889092

# Row 3
$ws.Range("G3").Value = 35.88321066666667
$ws.Range("H3").Value = 107.649632
$ws.Range("I3").Value = 0.08317795499144418
$ws.Range("J3").Value = 0.08448843719082051
$ws.Range("M3").Value = 53.72012833333334
$ws.Range("N3").Value = 161.160385
$ws.Range("O3").Value = 0.64508930810036
$ws.Range("P3").Value = 0.6636236722513976
$ws.Range("Q3").Value = 1927.650682025369
$ws.Range("R3").Value = 17348.85613822832
$ws.Range("S3").Value = 0.05365720943463361
$ws.Range("T3").Value = 0.05606852695135386

# Row 4
$ws.Range("G4").Value = 35.88321066666667
$ws.Range("H4").Value = 107.649632
$ws.Range("I4").Value = 0.08317795499144418
$ws.Range("J4").Value = 0.08448843719082051
$ws.Range("M4").Value = 0.2338196666666666
$ws.Range("N4").Value = 0.7014589999999999
$ws.Range("O4").Value = 0.002807784934062862
$ws.Range("P4").Value = 0.002888456722871399
$ws.Range("Q4").Value = 8.390200357009778
$ws.Range("R4").Value = 75.511803213088
$ws.Range("S4").Value = 0.0002335458088711357
$ws.Range("T4").Value = 0.0002440411944087235

# Row 5
$ws.Range("G5").Value = 35.88321066666667
$ws.Range("H5").Value = 107.649632
$ws.Range("I5").Value = 0.08317795499144418
$ws.Range("J5").Value = 0.08448843719082051
$ws.Range("M5").Value = 6.9774105
$ws.Range("N5").Value = 13.954821
$ws.Range("O5").Value = 0.08378708412325747
$ws.Range("P5").Value = 0.05746294014891389
$ws.Range("Q5").Value = 250.371890879312
$ws.Range("R5").Value = 1502.231345275872
$ws.Range("S5").Value = 0.006969238312068657
$ws.Range("T5").Value = 0.00485495400957139

# Row 6
$ws.Range("G6").Value = 35.88321066666667
$ws.Range("H6").Value = 107.649632
$ws.Range("I6").Value = 0.08317795499144418
$ws.Range("J6").Value = 0.08448843719082051
$ws.Range("M6").Value = 0.02062566666666667
$ws.Range("N6").Value = 0.061877
$ws.Range("O6").Value = 0.0002476799190900789
$ws.Range("P6").Value = 0.0002547961272734595
$ws.Range("Q6").Value = 0.7401151421404445
$ws.Range("R6").Value = 6.661036279264001
$ws.Range("S6").Value = 0.00002060150916235912
$ws.Range("T6").Value = 0.00002152732659560799

# Row 7
$ws.Range("I7").Value = 0.03522729558434242
$ws.Range("J7").Value = 0.03578230735158529
$ws.Range("M7").Value = 22.323506
$ws.Range("N7").Value = 66.970518
$ws.Range("O7").Value = 0.2680681429232296
$ws.Range("P7").Value = 0.2757701347495435
$ws.Range("Q7").Value = 339.2538329136106
$ws.Range("R7").Value = 3053.284496222496
$ws.Range("S7").Value = 0.009443315707502358
$ws.Range("T7").Value = 0.009867691719996258

# Row 8
$ws.Range("I8").Value = 0.03522729558434242
$ws.Range("J8").Value = 0.03578230735158529
$ws.Range("M8").Value = 53.72012833333334
$ws.Range("N8").Value = 161.160385
$ws.Range("O8").Value = 0.64508930810036
$ws.Range("P8").Value = 0.6636236722513976
$ws.Range("Q8").Value = 816.3932422485245
$ws.Range("R8").Value = 7347.53918023672
$ws.Range("S8").Value = 0.02272475173475032
$ws.Range("T8").Value = 0.02374598620628722

# Row 9
$ws.Range("I9").Value = 0.03522729558434242
$ws.Range("J9").Value = 0.03578230735158529
$ws.Range("M9").Value = 0.2338196666666666
$ws.Range("N9").Value = 0.7014589999999999
$ws.Range("O9").Value = 0.002807784934062862
$ws.Range("P9").Value = 0.002888456722871399
$ws.Range("Q9").Value = 3.553394261960888
$ws.Range("R9").Value = 31.98054835764799
$ws.Range("S9").Value = 0.00009891066980949582
$ws.Range("T9").Value = 0.0001033556462295372

# Row 10
$ws.Range("I10").Value = 0.03522729558434242
$ws.Range("J10").Value = 0.03578230735158529
$ws.Range("M10").Value = 6.9774105
$ws.Range("N10").Value = 13.954821
$ws.Range("O10").Value = 0.08378708412325747
$ws.Range("P10").Value = 0.05746294014891389
$ws.Range("Q10").Value = 106.036805147752
$ws.Range("R10").Value = 636.220830886512
$ws.Range("S10").Value = 0.002951592378560155
$ws.Range("T10").Value = 0.002056156585734187

# Row 11
$ws.Range("I11").Value = 0.03522729558434242
$ws.Range("J11").Value = 0.03578230735158529
$ws.Range("M11").Value = 0.02062566666666667
$ws.Range("N11").Value = 0.061877
$ws.Range("O11").Value = 0.0002476799190900789
$ws.Range("P11").Value = 0.0002547961272734595
$ws.Range("Q11").Value = 0.3134515014382222
$ws.Range("R11").Value = 2.821063512944
$ws.Range("S11").Value = 0.000008725093720092226
$ws.Range("T11").Value = 0.000009117193338092571

# Row 12
$ws.Range("G12").Value = 177.70077
$ws.Range("H12").Value = 533.10231
$ws.Range("I12").Value = 0.4119137160358794
$ws.Range("J12").Value = 0.4184034835782469
$ws.Range("M12").Value = 22.323506
$ws.Range("N12").Value = 66.970518
$ws.Range("O12").Value = 0.2680681429232296
$ws.Range("P12").Value = 0.2757701347495435
$ws.Range("Q12").Value = 3966.90420529962
$ws.Range("R12").Value = 35702.13784769658
$ws.Range("S12").Value = 0.1104209449023447
$ws.Range("T12").Value = 0.1153831850460516

# Row 13
$ws.Range("G13").Value = 177.70077
$ws.Range("H13").Value = 533.10231
$ws.Range("I13").Value = 0.4119137160358794
$ws.Range("J13").Value = 0.4184034835782469
$ws.Range("M13").Value = 53.72012833333334
$ws.Range("N13").Value = 161.160385
$ws.Range("O13").Value = 0.64508930810036
$ws.Range("P13").Value = 0.6636236722513976
$ws.Range("Q13").Value = 9546.108169332152
$ws.Range("R13").Value = 85914.97352398936
$ws.Range("S13").Value = 0.2657211340746336
$ws.Range("T13").Value = 0.2776624562549735

# Row 14
$ws.Range("G14").Value = 177.70077
$ws.Range("H14").Value = 533.10231
$ws.Range("I14").Value = 0.4119137160358794
$ws.Range("J14").Value = 0.4184034835782469
$ws.Range("M14").Value = 0.2338196666666666
$ws.Range("N14").Value = 0.7014589999999999
$ws.Range("O14").Value = 0.002807784934062862
$ws.Range("P14").Value = 0.002888456722871399
$ws.Range("Q14").Value = 41.54993480781
$ws.Range("R14").Value = 373.9494132702899
$ws.Range("S14").Value = 0.00115656512601939
$ws.Range("T14").Value = 0.0012085403550144

# Row 15
$ws.Range("G15").Value = 177.70077
$ws.Range("H15").Value = 533.10231
$ws.Range("I15").Value = 0.4119137160358794
$ws.Range("J15").Value = 0.4184034835782469
$ws.Range("M15").Value = 6.9774105
$ws.Range("N15").Value = 13.954821
$ws.Range("O15").Value = 0.08378708412325747
$ws.Range("P15").Value = 0.05746294014891389
$ws.Range("Q15").Value = 1239.891218456085
$ws.Range("R15").Value = 7439.34731073651
$ws.Range("S15").Value = 0.03451304917702182
$ws.Range("T15").Value = 0.02404269433495388

# Row 16
$ws.Range("G16").Value = 177.70077
$ws.Range("H16").Value = 533.10231
$ws.Range("I16").Value = 0.4119137160358794
$ws.Range("J16").Value = 0.4184034835782469
$ws.Range("M16").Value = 0.02062566666666667
$ws.Range("N16").Value = 0.061877
$ws.Range("O16").Value = 0.0002476799190900789
$ws.Range("P16").Value = 0.0002547961272734595
$ws.Range("Q16").Value = 3.66519684843
$ws.Range("R16").Value = 32.98677163587
$ws.Range("S16").Value = 0.0001020227558598604
$ws.Range("T16").Value = 0.0001066075872534618

# Row 17
$ws.Range("G17").Value = 20.074196
$ws.Range("H17").Value = 40.148392
$ws.Range("I17").Value = 0.04653236263856699
$ws.Range("J17").Value = 0.0315103250497358
$ws.Range("M17").Value = 22.323506
$ws.Range("N17").Value = 66.970518
$ws.Range("O17").Value = 0.2680681429232296
$ws.Range("P17").Value = 0.2757701347495435
$ws.Range("Q17").Value = 448.126434851176
$ws.Range("R17").Value = 2688.758609107056
$ws.Range("S17").Value = 0.01247384403835092
$ws.Range("T17").Value = 0.008689606584967558

# Row 18
$ws.Range("G18").Value = 20.074196
$ws.Range("H18").Value = 40.148392
$ws.Range("I18").Value = 0.04653236263856699
$ws.Range("J18").Value = 0.0315103250497358
$ws.Range("M18").Value = 53.72012833333334
$ws.Range("N18").Value = 161.160385
$ws.Range("O18").Value = 0.64508930810036
$ws.Range("P18").Value = 0.6636236722513976
$ws.Range("Q18").Value = 1078.388385308487
$ws.Range("R18").Value = 6470.330311850921
$ws.Range("S18").Value = 0.03001752961878822
$ws.Range("T18").Value = 0.02091099762334087

# Row 19
$ws.Range("G19").Value = 20.074196
$ws.Range("H19").Value = 40.148392
$ws.Range("I19").Value = 0.04653236263856699
$ws.Range("J19").Value = 0.0315103250497358
$ws.Range("M19").Value = 0.2338196666666666
$ws.Range("N19").Value = 0.7014589999999999
$ws.Range("O19").Value = 0.002807784934062862
$ws.Range("P19").Value = 0.002888456722871399
$ws.Range("Q19").Value = 4.693741817321333
$ws.Range("R19").Value = 28.162450903928
$ws.Range("S19").Value = 0.000130652866762918
$ws.Range("T19").Value = 0.00009101621022977242

# Row 20
$ws.Range("G20").Value = 20.074196
$ws.Range("H20").Value = 40.148392
$ws.Range("I20").Value = 0.04653236263856699
$ws.Range("J20").Value = 0.0315103250497358
$ws.Range("M20").Value = 6.9774105
$ws.Range("N20").Value = 13.954821
$ws.Range("O20").Value = 0.08378708412325747
$ws.Range("P20").Value = 0.05746294014891389
$ws.Range("Q20").Value = 140.065905949458
$ws.Range("R20").Value = 560.263623797832
$ws.Range("S20").Value = 0.003898810982851535
$ws.Range("T20").Value = 0.00181067592240579

# Row 21
$ws.Range("G21").Value = 20.074196
$ws.Range("H21").Value = 40.148392
$ws.Range("I21").Value = 0.04653236263856699
$ws.Range("J21").Value = 0.0315103250497358
$ws.Range("M21").Value = 0.02062566666666667
$ws.Range("N21").Value = 0.061877
$ws.Range("O21").Value = 0.0002476799190900789
$ws.Range("P21").Value = 0.0002547961272734595
$ws.Range("Q21").Value = 0.4140436752973334
$ws.Range("R21").Value = 2.484262051784
$ws.Range("S21").Value = 0.00001152513181339048
$ws.Range("T21").Value = 0.00000802870879180056

# Row 22
$ws.Range("G22").Value = 182.547562
$ws.Range("H22").Value = 547.642686
$ws.Range("I22").Value = 0.423148670749767
$ws.Range("J22").Value = 0.4298154468296114
$ws.Range("M22").Value = 22.323506
$ws.Range("N22").Value = 66.970518
$ws.Range("O22").Value = 0.2680681429232296
$ws.Range("P22").Value = 0.2757701347495435
$ws.Range("Q22").Value = 4075.101595592372
$ws.Range("R22").Value = 36675.91436033135
$ws.Range("S22").Value = 0.1134326783483232
$ws.Range("T22").Value = 0.1185302636896372

# Row 23
$ws.Range("G23").Value = 182.547562
$ws.Range("H23").Value = 547.642686
$ws.Range("I23").Value = 0.423148670749767
$ws.Range("J23").Value = 0.4298154468296114
$ws.Range("M23").Value = 53.72012833333334
$ws.Range("N23").Value = 161.160385
$ws.Range("O23").Value = 0.64508930810036
$ws.Range("P23").Value = 0.6636236722513976
$ws.Range("Q23").Value = 9806.478457577125
$ws.Range("R23").Value = 88258.30611819413
$ws.Range("S23").Value = 0.2729686832375542
$ws.Range("T23").Value = 0.2852357052154421

# Row 24
$ws.Range("G24").Value = 182.547562
$ws.Range("H24").Value = 547.642686
$ws.Range("I24").Value = 0.423148670749767
$ws.Range("J24").Value = 0.4298154468296114
$ws.Range("M24").Value = 0.2338196666666666
$ws.Range("N24").Value = 0.7014589999999999
$ws.Range("O24").Value = 0.002807784934062862
$ws.Range("P24").Value = 0.002888456722871399
$ws.Range("Q24").Value = 42.68321009765266
$ws.Range("R24").Value = 384.148890878874
$ws.Range("S24").Value = 0.001188110462599922
$ws.Range("T24").Value = 0.001241503316988966

# Row 25
$ws.Range("G25").Value = 182.547562
$ws.Range("H25").Value = 547.642686
$ws.Range("I25").Value = 0.423148670749767
$ws.Range("J25").Value = 0.4298154468296114
$ws.Range("M25").Value = 6.9774105
$ws.Range("N25").Value = 13.954821
$ws.Range("O25").Value = 0.08378708412325747
$ws.Range("P25").Value = 0.05746294014891389
$ws.Range("Q25").Value = 1273.709275848201
$ws.Range("R25").Value = 7642.255655089207
$ws.Range("S25").Value = 0.03545439327275531
$ws.Range("T25").Value = 0.02469845929624864

# Row 26
$ws.Range("G26").Value = 182.547562
$ws.Range("H26").Value = 547.642686
$ws.Range("I26").Value = 0.423148670749767
$ws.Range("J26").Value = 0.4298154468296114
$ws.Range("M26").Value = 0.02062566666666667
$ws.Range("N26").Value = 0.061877
$ws.Range("O26").Value = 0.0002476799190900789
$ws.Range("P26").Value = 0.0002547961272734595
$ws.Range("Q26").Value = 3.765165164624667
$ws.Range("R26").Value = 33.886486481622
$ws.Range("S26").Value = 0.0001048054285343767
$ws.Range("T26").Value = 0.0001095153112944965

Write-Host "Applied updates"